$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rtn4"
$ws.Range("C2").Value = "Tnfrsf19"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 57.045267
$ws.Range("H2").Value = 171.135801
$ws.Range("I2").Value = 0.2489699905037019
$ws.Range("J2").Value = 0.2489699905037019
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.08774700000000001
$ws.Range("N2").Value = 0.263241
$ws.Range("O2").Value = 0.1364824129982362
$ws.Range("P2").Value = 0.1364824129982362
$ws.Range("Q2").Value = 5.005551043449
$ws.Range("R2").Value = 45.04995939104101
$ws.Range("S2").Value = 0.03398002506809318
$ws.Range("T2").Value = 0.03398002506809317

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rtn4"
$ws.Range("C3").Value = "Tnfrsf19"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 57.045267
$ws.Range("H3").Value = 171.135801
$ws.Range("I3").Value = 0.2489699905037019
$ws.Range("J3").Value = 0.2489699905037019
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.555171
$ws.Range("N3").Value = 1.665513
$ws.Range("O3").Value = 0.8635175870017638
$ws.Range("P3").Value = 0.8635175870017638
$ws.Range("Q3").Value = 31.669877925657
$ws.Range("R3").Value = 285.028901330913
$ws.Range("S3").Value = 0.2149899654356087
$ws.Range("T3").Value = 0.2149899654356087

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rtn4"
$ws.Range("C4").Value = "Tnfrsf19"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 54.09018966666667
$ws.Range("H4").Value = 162.270569
$ws.Range("I4").Value = 0.2360727666969011
$ws.Range("J4").Value = 0.2360727666969011
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.08774700000000001
$ws.Range("N4").Value = 0.263241
$ws.Range("O4").Value = 0.1364824129982362
$ws.Range("P4").Value = 0.1364824129982362
$ws.Range("Q4").Value = 4.746251872681
$ws.Range("R4").Value = 42.716266854129
$ws.Range("S4").Value = 0.03221978084196271
$ws.Range("T4").Value = 0.0322197808419627

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rtn4"
$ws.Range("C5").Value = "Tnfrsf19"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 54.09018966666667
$ws.Range("H5").Value = 162.270569
$ws.Range("I5").Value = 0.2360727666969011
$ws.Range("J5").Value = 0.2360727666969011
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.555171
$ws.Range("N5").Value = 1.665513
$ws.Range("O5").Value = 0.8635175870017638
$ws.Range("P5").Value = 0.8635175870017638
$ws.Range("Q5").Value = 30.029304687433
$ws.Range("R5").Value = 270.263742186897
$ws.Range("S5").Value = 0.2038529858549384
$ws.Range("T5").Value = 0.2038529858549384

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Rtn4"
$ws.Range("C6").Value = "Tnfrsf19"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 67.444722
$ws.Range("H6").Value = 202.334166
$ws.Range("I6").Value = 0.2943576685488177
$ws.Range("J6").Value = 0.2943576685488177
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.08774700000000001
$ws.Range("N6").Value = 0.263241
$ws.Range("O6").Value = 0.1364824129982362
$ws.Range("P6").Value = 0.1364824129982362
$ws.Range("Q6").Value = 5.918072021334
$ws.Range("R6").Value = 53.262648192006
$ws.Range("S6").Value = 0.04017464488807766
$ws.Range("T6").Value = 0.04017464488807765

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Rtn4"
$ws.Range("C7").Value = "Tnfrsf19"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 67.444722
$ws.Range("H7").Value = 202.334166
$ws.Range("I7").Value = 0.2943576685488177
$ws.Range("J7").Value = 0.2943576685488177
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.555171
$ws.Range("N7").Value = 1.665513
$ws.Range("O7").Value = 0.8635175870017638
$ws.Range("P7").Value = 0.8635175870017638
$ws.Range("Q7").Value = 37.443353757462
$ws.Range("R7").Value = 336.990183817158
$ws.Range("S7").Value = 0.2541830236607401
$ws.Range("T7").Value = 0.2541830236607401

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Rtn4"
$ws.Range("C8").Value = "Tnfrsf19"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 50.54489333333333
$ws.Range("H8").Value = 151.63468
$ws.Range("I8").Value = 0.2205995742505793
$ws.Range("J8").Value = 0.2205995742505793
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.08774700000000001
$ws.Range("N8").Value = 0.263241
$ws.Range("O8").Value = 0.1364824129982362
$ws.Range("P8").Value = 0.1364824129982362
$ws.Range("Q8").Value = 4.43516275532
$ws.Range("R8").Value = 39.91646479788
$ws.Range("S8").Value = 0.03010796220010264
$ws.Range("T8").Value = 0.03010796220010263

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Rtn4"
$ws.Range("C9").Value = "Tnfrsf19"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 50.54489333333333
$ws.Range("H9").Value = 151.63468
$ws.Range("I9").Value = 0.2205995742505793
$ws.Range("J9").Value = 0.2205995742505793
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.555171
$ws.Range("N9").Value = 1.665513
$ws.Range("O9").Value = 0.8635175870017638
$ws.Range("P9").Value = 0.8635175870017638
$ws.Range("Q9").Value = 28.06105897676
$ws.Range("R9").Value = 252.54953079084
$ws.Range("S9").Value = 0.1904916120504767
$ws.Range("T9").Value = 0.1904916120504767
